$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared-string runs) ---
$ws.Range("A8").Value = "Volume 30   Number  7"
$ws.Range("C9").Value = "Report Covering the Week  2/13/2023  Through  2/19/2023"

# --- Numeric cell updates ---
$ws.Range("C14").Value = 5
$ws.Range("D14").Value = 8
$ws.Range("E14").Value = -37.5
$ws.Range("F14").Value = 25
$ws.Range("G14").Value = 33
$ws.Range("H14").Value = -24.242424242424
$ws.Range("I14").Value = 49
$ws.Range("J14").Value = 56
$ws.Range("K14").Value = -12.5
$ws.Range("L14").Value = -2
$ws.Range("M14").Value = -9.259259259259
$ws.Range("N14").Value = -82.246376811594
$ws.Range("C15").Value = 30
$ws.Range("D15").Value = 41
$ws.Range("E15").Value = -26.829268292682
$ws.Range("F15").Value = 107
$ws.Range("G15").Value = 128
$ws.Range("H15").Value = -16.40625
$ws.Range("I15").Value = 210
$ws.Range("J15").Value = 225
$ws.Range("K15").Value = -6.666666666666
$ws.Range("L15").Value = 26.506024096385
$ws.Range("M15").Value = 36.363636363636
$ws.Range("N15").Value = -48.780487804878
$ws.Range("C16").Value = 311
$ws.Range("D16").Value = 367
$ws.Range("E16").Value = -15.258855585831
$ws.Range("F16").Value = 1164
$ws.Range("G16").Value = 1205
$ws.Range("H16").Value = -3.402489626556
$ws.Range("I16").Value = 2110
$ws.Range("J16").Value = 2119
$ws.Range("K16").Value = -0.424728645587
$ws.Range("L16").Value = 43.245078071962
$ws.Range("M16").Value = -16.036609629924
$ws.Range("N16").Value = -82.180559074402
$ws.Range("C17").Value = 483
$ws.Range("D17").Value = 450
$ws.Range("E17").Value = 7.333333333333
$ws.Range("F17").Value = 1819
$ws.Range("G17").Value = 1713
$ws.Range("H17").Value = 6.187974314068
$ws.Range("I17").Value = 3348
$ws.Range("J17").Value = 3002
$ws.Range("K17").Value = 11.525649566955
$ws.Range("L17").Value = 37.721102426984
$ws.Range("M17").Value = 70.816326530612
$ws.Range("N17").Value = -28.856778580535
$ws.Range("C18").Value = 267
$ws.Range("D18").Value = 313
$ws.Range("E18").Value = -14.696485623003
$ws.Range("F18").Value = 1078
$ws.Range("G18").Value = 1167
$ws.Range("H18").Value = -7.626392459297
$ws.Range("I18").Value = 2007
$ws.Range("J18").Value = 2044
$ws.Range("K18").Value = -1.810176125244
$ws.Range("L18").Value = 21.41560798548
$ws.Range("M18").Value = -20.70327933623
$ws.Range("N18").Value = -85.58810857389
$ws.Range("C19").Value = 849
$ws.Range("D19").Value = 963
$ws.Range("E19").Value = -11.838006230529
$ws.Range("F19").Value = 3553
$ws.Range("G19").Value = 3706
$ws.Range("H19").Value = -4.128440366972
$ws.Range("I19").Value = 6405
$ws.Range("J19").Value = 6712
$ws.Range("K19").Value = -4.57389749702
$ws.Range("L19").Value = 62.604722010662
$ws.Range("M19").Value = 35.987261146496
$ws.Range("N19").Value = -39.655172413793
$ws.Range("C20").Value = 277
$ws.Range("D20").Value = 279
$ws.Range("E20").Value = -0.716845878136
$ws.Range("F20").Value = 1128
$ws.Range("G20").Value = 1040
$ws.Range("H20").Value = 8.461538461538
$ws.Range("I20").Value = 1980
$ws.Range("J20").Value = 1910
$ws.Range("K20").Value = 3.664921465968
$ws.Range("L20").Value = 102.453987730061
$ws.Range("M20").Value = 53.132250580046
$ws.Range("N20").Value = -87.648930197741
$ws.Range("C21").Value = 2222
$ws.Range("D21").Value = 2421
$ws.Range("E21").Value = -8.219743907476
$ws.Range("F21").Value = 8874
$ws.Range("G21").Value = 8992
$ws.Range("H21").Value = -1.312277580071
$ws.Range("I21").Value = 16109
$ws.Range("J21").Value = 16068
$ws.Range("K21").Value = 0.255165546427
$ws.Range("L21").Value = 50.692235734331
$ws.Range("M21").Value = 21.899356791524
$ws.Range("N21").Value = -72.131686388485
$ws.Range("C22").Value = 43
$ws.Range("D22").Value = 48
$ws.Range("E22").Value = -10.416666666666
$ws.Range("F22").Value = 156
$ws.Range("G22").Value = 178
$ws.Range("H22").Value = -12.359550561797
$ws.Range("I22").Value = 259
$ws.Range("J22").Value = 323
$ws.Range("K22").Value = -19.814241486068
$ws.Range("L22").Value = 32.820512820512
$ws.Range("M22").Value = -8.480565371024
$ws.Range("C23").Value = 106
$ws.Range("D23").Value = 118
$ws.Range("E23").Value = -10.169491525423
$ws.Range("F23").Value = 437
$ws.Range("G23").Value = 443
$ws.Range("H23").Value = -1.354401805869
$ws.Range("I23").Value = 797
$ws.Range("J23").Value = 766
$ws.Range("K23").Value = 4.046997389033
$ws.Range("L23").Value = 23.374613003096
$ws.Range("M23").Value = 60.04016064257
$ws.Range("C24").Value = 1951
$ws.Range("D24").Value = 2094
$ws.Range("E24").Value = -6.829035339063
$ws.Range("F24").Value = 8121
$ws.Range("G24").Value = 7888
$ws.Range("H24").Value = 2.953853955375
$ws.Range("I24").Value = 14285
$ws.Range("J24").Value = 13323
$ws.Range("K24").Value = 7.22059596187
$ws.Range("L24").Value = 39.542834814887
$ws.Range("M24").Value = 43.021625951141
$ws.Range("C25").Value = 745
$ws.Range("D25").Value = 726
$ws.Range("E25").Value = 2.617079889807
$ws.Range("F25").Value = 3005
$ws.Range("G25").Value = 2917
$ws.Range("H25").Value = 3.016798080219
$ws.Range("I25").Value = 5348
$ws.Range("J25").Value = 4916
$ws.Range("K25").Value = 8.787632221318
$ws.Range("L25").Value = 43.763440860215
$ws.Range("M25").Value = -1.655020228025
$ws.Range("C26").Value = 50
$ws.Range("D26").Value = 59
$ws.Range("E26").Value = -15.254237288135
$ws.Range("F26").Value = 183
$ws.Range("G26").Value = 207
$ws.Range("H26").Value = -11.59420289855
$ws.Range("I26").Value = 335
$ws.Range("J26").Value = 343
$ws.Range("K26").Value = -2.332361516034
$ws.Range("L26").Value = 16.724738675958
$ws.Range("C27").Value = 89
$ws.Range("D27").Value = 106
$ws.Range("E27").Value = -16.037735849056
$ws.Range("F27").Value = 359
$ws.Range("G27").Value = 364
$ws.Range("H27").Value = -1.373626373626
$ws.Range("I27").Value = 653
$ws.Range("J27").Value = 591
$ws.Range("K27").Value = 10.490693739424
$ws.Range("L27").Value = 34.086242299794
$ws.Range("C28").Value = 16
$ws.Range("D28").Value = 16
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 80
$ws.Range("G28").Value = 100
$ws.Range("H28").Value = -20
$ws.Range("I28").Value = 149
$ws.Range("J28").Value = 173
$ws.Range("K28").Value = -13.872832369942
$ws.Range("L28").Value = 13.740458015267
$ws.Range("M28").Value = -4.487179487179
$ws.Range("N28").Value = -80.186170212766
$ws.Range("C29").Value = 14
$ws.Range("D29").Value = 15
$ws.Range("E29").Value = -6.666666666666
$ws.Range("F29").Value = 66
$ws.Range("G29").Value = 88
$ws.Range("H29").Value = -25
$ws.Range("I29").Value = 123
$ws.Range("J29").Value = 154
$ws.Range("K29").Value = -20.12987012987
$ws.Range("L29").Value = 2.5
$ws.Range("M29").Value = -11.510791366906
$ws.Range("N29").Value = -82.327586206896
$ws.Range("D30").Value = 27
$ws.Range("E30").Value = -85.185185185185
$ws.Range("F30").Value = 22
$ws.Range("G30").Value = 68
$ws.Range("H30").Value = -67.647058823529
$ws.Range("I30").Value = 45
$ws.Range("J30").Value = 96
$ws.Range("K30").Value = -53.125
$ws.Range("L30").Value = 36.363636363636
